# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Changes applied:
#  1. Metadata sheet: Version 1.0.0 -> 0.1.0, Status active -> draft,
#     Experimental (blank) -> false, Date updated, Description filled in.
#  2. "Include #0" sheet: collapses the old 3-column Property/Operation/Value
#     layout into a 2-column Concept/Description layout with the concepts
#     that belong to the "At risk for falls" style include, adding two new
#     concept rows, and keeps the trailing System URI row.
#  3. New "Include #1" sheet added (same 2-column Concept/Description shape)
#     for the "Risk of falls" include pointing at the local CodeSystem.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.1.0"
$meta.Range("B6").Value = "draft"

# "false" needs to land as literal text (matches the shared-string in the
# target, not a Boolean) - a leading apostrophe forces text entry, then we
# re-apply the plain body style so the quote-prefix marker doesn't linger.
$meta.Range("B7").Value = "'false"
$meta.Range("B6").Copy()
$meta.Range("B7").PasteSpecial(-4122)

$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"
$meta.Range("B11").Value = "Value set for nursing problems and diagnoses"

# ---------------------------------------------------------------------
# 2. "Include #0" sheet - rebuild as a 2-column Concept/Description table
# ---------------------------------------------------------------------
$inc0 = $wb.Worksheets.Item("Include #0")

# Rows 5 and 6 don't exist yet - seed them with the existing body-row
# style (s="2") by pasting row 4's format, so they come out with the
# right style even for the cells that end up blank.
$inc0.Range("A4:B4").Copy()
$inc0.Range("A5:B6").PasteSpecial(-4122)

# Header row (style s=1 already present on row 1).
$inc0.Range("A1").Value = "Concept"
$inc0.Range("B1").Value = "Description"

# Body rows - the SNOMED codes are text-looking digit strings in the
# source data, so force text entry with a leading apostrophe (otherwise
# they'd be auto-coerced to numbers).
$inc0.Range("A2").Value = "'129839007"
$inc0.Range("B2").Value = "At risk for falls"

$inc0.Range("A3").Value = "'300893006"
$inc0.Range("B3").Value = "Nutritional finding"

$inc0.Range("A4").Value = "'22253000"
$inc0.Range("B4").Value = "Pain"

$inc0.Range("A5").Value = ""
$inc0.Range("B5").Value = ""

$inc0.Range("A6").Value = "System URI"
$inc0.Range("B6").Value = "http://snomed.info/sct"

# Drop the now-unused third column (old "Value" column).
$inc0.Range("C1:C4").Clear()

# Strip the quote-prefix markers the apostrophe entry left on column A by
# re-pasting the plain body format (sourced from an untouched cell) back
# over those rows.
$inc0.Range("B2").Copy()
$inc0.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. New "Include #1" sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$inc1 = $wb.Worksheets.Add($null, $lastSheet)
$inc1.Name = "Include #1"

# Match column widths and the header/body styling from "Include #0".
$inc1.Range("A1").ColumnWidth = 29.869791666666668
$inc1.Range("B1").ColumnWidth = 49.869791666666664

$inc0.Range("A1:B1").Copy()
$inc1.Range("A1:B1").PasteSpecial(-4122)
$inc0.Range("A2:B2").Copy()
$inc1.Range("A2:B4").PasteSpecial(-4122)

$inc1.Range("A1").Value = "Concept"
$inc1.Range("B1").Value = "Description"

$inc1.Range("A2").Value = "risk-falls"
$inc1.Range("B2").Value = "Risk of falls"

$inc1.Range("A3").Value = ""
$inc1.Range("B3").Value = ""

$inc1.Range("A4").Value = "System URI"
$inc1.Range("B4").Value = "https://clinyqai.github.io/open-nursing-core-ig/CodeSystem/onc-observation-codes"
